$d = $word.ActiveDocument

$replacements = @(
    @("840×9=7560", "377×9=3393"),
    @("832×9=7488", "547×4=2188"),
    @("223×4=892", "391×8=3128"),
    @("858×9=7722", "517×3=1551"),
    @("105×9=945", "794×7=5558"),
    @("523×6=3138", "835×4=3340"),
    @("528×9=4752", "609×3=1827"),
    @("860×5=4300", "428×3=1284"),
    @("921×3=2763", "981×5=4905"),
    @("485×6=2910", "137×6=822"),
    @("383×6=2298", "780×4=3120"),
    @("144×4=576", "401×5=2005"),
    @("419×6=2514", "420×8=3360"),
    @("164×5=820", "371×2=742"),
    @("401×4=1604", "184×2=368"),
    @("975×4=3900", "359×3=1077"),
    @("248×2=496", "687×2=1374"),
    @("491×4=1964", "647×8=5176"),
    @("125×3=375", "143×7=1001"),
    @("703×6=4218", "689×6=4134"),
    @("242×3=726", "505×6=3030"),
    @("536×7=3752", "550×4=2200"),
    @("456×4=1824", "872×7=6104"),
    @("305×3=915", "704×8=5632"),
    @("149×3=447", "797×3=2391")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
